$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B3: strain1 -> MGP47
$ws.Range("B3").Value = "MGP47"

# New row 8: Header Format / TIME::VALUE_TYPE / Must be TIME::VALUE_TYPE
$ws.Range("A8").Value = "Header Format"
$ws.Range("B8").Value = "TIME::VALUE_TYPE"
$ws.Range("C8").Value = "Must be TIME::VALUE_TYPE"

# Formatting for the new row, matching the look of the existing rows
$ws.Range("A8").Font.Name = "Verdana"
$ws.Range("A8").Font.Size = 14
$ws.Range("A8").Font.Bold = $false
$ws.Range("A8").Font.Italic = $false

$ws.Range("B8").Font.Name = "Verdana"
$ws.Range("B8").Font.Size = 14
$ws.Range("B8").Font.Bold = $true
$ws.Range("B8").Font.Italic = $false

# Copy the formatting (font/border) of an existing description cell for C8
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Value = "Must be TIME::VALUE_TYPE"

# Update selection on sheet1 to match the new active cell
$ws.Select() | Out-Null
$ws.Range("B12").Select() | Out-Null

# Update the workbook window view geometry
$win = $excel.ActiveWindow
$win.Left = -20
$win.Top = -20
$win.Width = 24720
$win.Height = 16740
